$wb = $excel.ActiveWorkbook
$excel.ActiveWindow.Left = 11800
$excel.ActiveWindow.Top = 1760
